$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = 45233

for ($r = 2; $r -le 173; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
